# 3DES projetos atualizada a estrutura de pastas
# Adds an "AT1-Cont" column to FREQ (col O) and PONTOS (col E), and
# switches the active sheet/selection back to FREQ.

$wb = $excel.ActiveWorkbook

$freq = $wb.Worksheets.Item("FREQ")
$pontos = $wb.Worksheets.Item("PONTOS")

# --- FREQ sheet: fill column O with the new "AT1-Cont" attendance values ---
$freqValues = @{
    3  = "P"
    4  = "P"
    5  = "P"
    6  = "P"
    7  = "F"
    8  = "F"
    9  = "P"
    10 = "P"
    11 = "P"
    12 = "P"
    13 = "P"
    14 = "P"
    15 = "P"
    16 = "P"
    17 = "P"
    18 = "P"
    19 = "P"
    20 = "P"
}

foreach ($row in $freqValues.Keys) {
    $freq.Cells.Item($row, 15).Value = $freqValues[$row]
}

# --- PONTOS sheet: add the "AT1-Cont" header in column E ---
$pontos.Range("E1").Value = "AT1-Cont"

# --- Selection / active sheet bookkeeping ---
[void]$pontos.Range("E2").Select()
[void]$freq.Select()
[void]$freq.Range("O3").Select()
